$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force Text format first so numeric-looking
# strings (which Excel would otherwise auto-convert to Number) are stored
# verbatim as text, matching the source data exactly.
$ws.Range("D2").Value = "42.821.20"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "2.563.85"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.79"
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.59"
$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.76"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.46"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "2.958.82"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("E15").Value = "  +5.28%  "

$ws.Range("D16").Value = "2.610.98"
$ws.Range("E16").Value = "  +3.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "42.847.79"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  -1.22%  "

$ws.Range("D20").Value = "0.0₃0960"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.38"
$ws.Range("E21").Value = "  -2.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.39"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.89"
$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.04"
$ws.Range("E26").Value = "  +2.13%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.87"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.18"
$ws.Range("E30").Value = "  -2.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.11"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("E32").Value = "  -2.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0797"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.10"
$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.66"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.58"
$ws.Range("E38").Value = "  +12.88%  "

$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.76"
$ws.Range("E41").Value = "  +1.25%  "

$ws.Range("E42").Value = "  +7.75%  "

$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("D46").Value = "1.990.91"
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "2.810.54"
$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.34"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.77"
$ws.Range("E51").Value = "  -1.66%  "
